$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as text in the
# source data (e.g. "156.37", "3.00", "0.116"). Force these cells to the
# Text number format before writing so Excel keeps them as strings
# instead of silently converting them to numbers (which would also
# strip meaningful trailing zeros).
$priceCells = @("D2","D3","D5","D6","D7","D8","D10","D11","D12","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D28","D29","D30","D32","D33","D35","D36","D38","D40","D41","D42","D43","D44","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "42.553.45"
$ws.Range("E2").Value = "  +1.74%  "

# Row 3
$ws.Range("D3").Value = "2.290.93"
$ws.Range("E3").Value = "  +1.07%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "156.37"
$ws.Range("E5").Value = "  +15,518.85%  "

# Row 6
$ws.Range("D6").Value = "307.73"
$ws.Range("E6").Value = "  +0.86%  "

# Row 7
$ws.Range("D7").Value = "96.12"
$ws.Range("E7").Value = "  +4.47%  "

# Row 8
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").Value = "0.497"
$ws.Range("E10").Value = "  +2.86%  "

# Row 11
$ws.Range("D11").Value = "35.68"
$ws.Range("E11").Value = "  +9.25%  "

# Row 12
$ws.Range("D12").Value = "0.0807"
$ws.Range("E12").Value = "  +1.23%  "

# Row 13
$ws.Range("E13").Value = "  -1.61%  "

# Row 14
$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  +1.87%  "

# Row 15
$ws.Range("D15").Value = "2.649.83"
$ws.Range("E15").Value = "  +1.22%  "

# Row 16
$ws.Range("D16").Value = "14.55"
$ws.Range("E16").Value = "  +2.23%  "

# Row 17
$ws.Range("D17").Value = "2.306.25"
$ws.Range("E17").Value = "  +1.55%  "

# Row 18
$ws.Range("D18").Value = "0.801"
$ws.Range("E18").Value = "  +4.88%  "

# Row 19
$ws.Range("D19").Value = "42.455.35"
$ws.Range("E19").Value = "  +1.69%  "

# Row 20
$ws.Range("D20").Value = "12.71"
$ws.Range("E20").Value = "  +4.07%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0923"
$ws.Range("E21").Value = "  +1.82%  "

# Row 22
$ws.Range("D22").Value = "6.02"
$ws.Range("E22").Value = "  +1.82%  "

# Row 23
$ws.Range("D23").Value = "68.19"
$ws.Range("E23").Value = "  +1.97%  "

# Row 24
$ws.Range("D24").Value = "243.49"
$ws.Range("E24").Value = "  +0.91%  "

# Row 25
$ws.Range("D25").Value = "2.61"
$ws.Range("E25").Value = "  +0.46%  "

# Row 26
$ws.Range("D26").Value = "1.96"
$ws.Range("E26").Value = "  +2.14%  "

# Row 27
$ws.Range("E27").Value = "  -0.26%  "

# Row 28
$ws.Range("D28").Value = "24.19"
$ws.Range("E28").Value = "  +0.13%  "

# Row 29
$ws.Range("D29").Value = "36.61"
$ws.Range("E29").Value = "  +7.52%  "

# Row 30
$ws.Range("D30").Value = "9.69"
$ws.Range("E30").Value = "  +0.78%  "

# Row 31
$ws.Range("E31").Value = "  -8.58%  "

# Row 32
$ws.Range("D32").Value = "161.35"
$ws.Range("E32").Value = "  +2.32%  "

# Row 33
$ws.Range("D33").Value = "5.36"
$ws.Range("E33").Value = "  +3.91%  "

# Row 34
$ws.Range("E34").Value = "  +0.07%  "

# Row 35
$ws.Range("D35").Value = "0.0757"
$ws.Range("E35").Value = "  +1.30%  "

# Row 36
$ws.Range("D36").Value = "3.10"
$ws.Range("E36").Value = "  +2.50%  "

# Row 37
$ws.Range("E37").Value = "  +4.93%  "

# Row 38
$ws.Range("D38").Value = "17.26"
$ws.Range("E38").Value = "  +0.04%  "

# Row 40 (content swapped with row 41: ARBITRUM -> Stellar)
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "0.116"
$ws.Range("E40").Value = "  -0.26%  "

# Row 41 (content swapped with row 40: Stellar -> ARBITRUM)
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.84"
$ws.Range("E41").Value = "  +1.69%  "

# Row 42
$ws.Range("D42").Value = "4.21"
$ws.Range("E42").Value = "  +7.42%  "

# Row 43 (content swapped with row 44: EnergySwap -> Maker)
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.020.13"
$ws.Range("E43").Value = "  -2.11%  "

# Row 44 (content swapped with row 43: Maker -> EnergySwap)
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "19.69"
$ws.Range("E44").Value = "  +1.54%  "

# Row 45
$ws.Range("E45").Value = "  +11.22%  "

# Row 46
$ws.Range("E46").Value = "  +2.24%  "

# Row 47
$ws.Range("D47").Value = "10.27"
$ws.Range("E47").Value = "  -0.14%  "

# Row 48
$ws.Range("D48").Value = "3.00"
$ws.Range("E48").Value = "  +3.47%  "

# Row 49
$ws.Range("D49").Value = "53.53"
$ws.Range("E49").Value = "  +3.72%  "

# Row 50
$ws.Range("D50").Value = "1.54"
$ws.Range("E50").Value = "  +0.67%  "

# Row 51
$ws.Range("D51").Value = "73.56"
$ws.Range("E51").Value = "  +1.46%  "
